# Journal de bord — "Version 1.0 avec Documentation"
#
# Fixes a handful of typos/capitalisation issues in the "Branche"/"Type"/
# "Description" columns and appends two new journal rows (20 & 21) for the
# "Version 1.0" delivery and its documentation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("journal")

# --- 1. Text corrections in existing rows -----------------------------
# Row 5 used to read "Git + Github" / "Instalation git + création dépôt
# github" — fix capitalisation & spelling.
$ws.Range("H5").Value = "Git + GitHub"
$ws.Range("I5").Value = "Installation git + création dépôt GitHub"

# Row 17 ("Selection de maps aléatoire" -> add missing accent).
$ws.Range("I17").Value = "Sélection de maps aléatoire"

# Row 19 ("Efface l'écran au fure et a mesure" -> fix typo/accents).
$ws.Range("I19").Value = "Efface l'écran au fur et à mesure"

# --- 2. Append two new rows, copying the formatting of row 19 ---------
$ws.Range("A19:J19").Copy()
$ws.Range("A20:J21").PasteSpecial(-4122)  # xlPasteFormats

# Row 20: Version 1.0
$ws.Range("A20").Value = 43929
$ws.Range("B20").Value = 7
$ws.Range("G20").Value = "MA-20"
$ws.Range("H20").Value = "Version 1.0"
$ws.Range("I20").Value = "Rendu de la Version 1.0"

# Row 21: Documentation
$ws.Range("A21").Value = 43929
$ws.Range("B21").Value = 7
$ws.Range("G21").Value = "ICT-431"
$ws.Range("H21").Value = "Documentation"
$ws.Range("I21").Value = "Rendu de la Documentation"

# Fill the duration formula down across the new rows too, which also
# regroups F12:F21 into a single shared formula like the rest of the sheet.
$ws.Range("F12:F21").Formula = '=IF(AND(C12<>"",D12<>""),D12-C12-E12,"")'

# --- 3. Selection / view bookkeeping -----------------------------------
$ws.Range("I19").Select() | Out-Null
